# Trade #21 closed at 2026-02-17 12:37:03 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---- Summary sheet ----
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.57   # Current Capital
$summary.Range("B4").Value = 0.57      # Total P&L $
$summary.Range("B5").Value = 0.54      # Total P&L %
$summary.Range("B6").Value = 21        # Total Trades
$summary.Range("B7").Value = 7         # Winning Trades
$summary.Range("B9").Value = 33.33     # Win Rate %

# ---- Strategy Status sheet (MarketMaking row) ----
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.57     # Capital
$status.Range("D4").Value = 21         # Trades
$status.Range("E4").Value = 0.57       # P&L $
$status.Range("F4").Value = 0.57       # P&L %
$status.Range("G4").Value = 33.33      # Win Rate %

# ---- All Trades sheet (trade #21 row, row 22) ----
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("G22").Value = 0.59          # Exit Price
$allTrades.Range("H22").Value = "CLOSED"      # Status
$allTrades.Range("I22").Value = 5800          # P&L %
$allTrades.Range("J22").Value = 0.58          # P&L $
$allTrades.Range("K22").Value = 100.57        # Capital After
$allTrades.Range("P22").Value = "early_exit"  # Exit Reason
$allTrades.Range("Q22").Value = 7.62          # Duration (min)

# ---- MarketMaking sheet (trade #21 row, row 22) ----
$marketMaking = $wb.Worksheets.Item("MarketMaking")
$marketMaking.Range("G22").Value = 0.59          # Exit Price
$marketMaking.Range("H22").Value = "CLOSED"      # Status
$marketMaking.Range("I22").Value = 5800          # P&L %
$marketMaking.Range("J22").Value = 0.58          # P&L $
$marketMaking.Range("K22").Value = 100.57        # Capital After
$marketMaking.Range("P22").Value = "early_exit"  # Exit Reason
$marketMaking.Range("Q22").Value = 7.62          # Duration (min)
